$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny rounding difference in A23
$ws.Range("A23").Value = 45874.91690625

# Add new row 24 data
$ws.Range("A24").Value = 45874.95854735681
$ws.Range("B24").Value = 2025
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 13.72
$ws.Range("E24").Value = 91.44
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 4.33
$ws.Range("H24").Value = "ESE"
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = "23:00:18"

# Apply same number format as A23 (date/time) to A24
$ws.Range("A24").NumberFormat = $ws.Range("A23").NumberFormat
